# Applies the "Adjustments to DSGE class" edit:
# - On sheet "Iterações": C1 becomes a hard-coded 20 (formula removed),
#   C2 becomes 2.1 (was 2.15), C3 keeps its formula and recalculates.
# - Selections / zoom updated to match the saved view state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Iterações")

# --- Data changes on "Iterações" ---
# C1 previously held formula =6+7 (cached 13); now a literal value.
$ws2.Range("C1").Value = 20
# C2 previously 2.15; now 2.1
$ws2.Range("C2").Value = 2.1

# Recalculate so C3 (=C1*3600*C2) reflects the new inputs (151200)
$excel.Calculate()

# --- View state changes ---
$ws1.Activate()
$ws1.Application.ActiveWindow.Zoom = 140
$ws1.Range("D19").Select()

$ws2.Activate()
$ws2.Range("D3").Select()

$ws1.Activate()
